# The commit removes the "Classification: Controlled" text-box shapes
# (an mc:AlternateContent drawing/VML fallback pair) that were anchored
# in each of the document's three footers, leaving each footer as a
# single empty paragraph styled "Footer".
#
# wdHeaderFooterIndex values used by Section.Footers():
#   1 = wdHeaderFooterPrimary   (footer2.xml - default footer)
#   2 = wdHeaderFooterFirstPage (footer3.xml - first-page footer)
#   3 = wdHeaderFooterEvenPages (footer1.xml - even-page footer)

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {
    for ($i = 1; $i -le 3; $i++) {
        $ftr = $sec.Footers($i)
        # Remove every shape (the classification text box) anchored in
        # this footer, regardless of whether the feature toggle makes
        # the footer "Exist" for this particular section/page type -
        # the underlying footer part still carries the shape.
        for ($j = $ftr.Shapes.Count; $j -ge 1; $j--) {
            $ftr.Shapes($j).Delete()
        }
    }
}
